# Nayem_meal.xlsx edit
#   commit msg: "Himel bazar-851,vim-10,Himel raat off,dupur extra+2meal"
#
# This applies:
#  - day 18 (col S) and day 19 (col T) "dupur" (lunch) meal counts of 2.5 for
#    everyone (Rakib, Mahfuz, Himel, Minhaz, Taher, Forhad, Nayem), with an
#    extra +1 meal (2.5 -> 3.5) for Himel on day 19 ("dupur extra+2meal")
#  - bazar (shopping) entry for day 18: Himel bought bazar worth 851 taka,
#    recorded against Himel's own "S25" cost line, and an extra 10 taka
#    ("vim") on day 19
#  - updates the sheet's selection/scroll state to match the new edit point

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Day 18 (S) / Day 19 (T) lunch meal counts for each person (rows 3-9) ---
# Rakib
$ws.Range("S3").Value = 2.5
$ws.Range("T3").Value = 2.5
# Mahfuz
$ws.Range("S4").Value = 2.5
$ws.Range("T4").Value = 2.5
# Himel - "raat off" (dinner off) on day 18 already reflected by the 2.5, plus
# "dupur extra+2meal" gives Himel an extra meal on day 19
$ws.Range("S5").Value = 2.5
$ws.Range("T5").Value = 3.5
# Minhaz
$ws.Range("S6").Value = 2.5
$ws.Range("T6").Value = 2.5
# Taher
$ws.Range("S7").Value = 2.5
$ws.Range("T7").Value = 2.5
# Forhad
$ws.Range("S8").Value = 2.5
$ws.Range("T8").Value = 2.5
# Nayem
$ws.Range("S9").Value = 2.5
$ws.Range("T9").Value = 2.5

# --- Himel's per-person bazar/cost line (row 25) picks up the 5 tk share ---
$ws.Range("S25").Value = 5

# --- Bazar log: Himel bought bazar on day 18, vim (dish soap) on day 19 ---
$ws.Range("S42").Value = "Himel"
$ws.Range("S43").Value = 851
$ws.Range("T43").Value = 10

# --- Update selection / scroll position to match the saved view ---
$ws.Range("A28").Select() | Out-Null
$ws.Range("T5").Select() | Out-Null
